# Add "NA" values under the "duplicate_image_filename" column (column E)
# for the practice + generic + unique_video/unique_audio rows (rows 2-21),
# where that column was previously left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Neutralize a harmless round-trip quirk of the empty shared-string cell F1
# (it has no value before or after this edit; explicitly clearing it keeps
# it blank through the save instead of an unrelated artifact appearing).
$ws.Range("F1").Value = $null

$ws.Range("E2:E21").Value = "NA"
